$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.589.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "'1.690.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'314.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.3889"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("D8").Value = "'0.4038"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'1.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'1.002"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'53.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").Value = "'0.08762"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "'25.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.38%  "
$ws.Range("D14").Value = "'7.511"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "'0.00001354"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "'7.955"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "'1.690.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'98.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "'0.07096"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'19.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'7.281"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.79%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'14.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'24.588.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'2.974"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.45%  "
$ws.Range("D26").Value = "'2.354"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'22.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'161.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").Value = "'8.818"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.40%  "
$ws.Range("D30").Value = "'137.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'5.222"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "'1.876.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'0.08820"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "'7.393"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.72%  "
$ws.Range("D35").Value = "'1.037"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("D37").Value = "'0.2755"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'0.02925"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.24%  "
$ws.Range("D39").Value = "'10.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.33%  "
$ws.Range("D40").Value = "'14.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").Value = "'0.09132"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "'0.7924"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").Value = "'1.457"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'16.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("D45").Value = "'0.7220"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "'2.597"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'4.201"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'1.349"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'138.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'91.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "

Write-Host "Updated crypto prices"
